$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12703
$ws1.Range("F9").Value = 3
$ws1.Range("F10").Value = 12603
$ws1.Range("F11").Value = 256
$ws1.Range("F12").Value = 7
$ws1.Range("F13").Value = 4936
$ws1.Range("F14").Value = 5619
$ws1.Range("F15").Value = 171
$ws1.Range("F16").Value = 80
$ws1.Range("F19").Value = 974
$ws1.Range("F23").Value = 179

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12703
$ws4.Range("F10").Value = 3
$ws4.Range("F11").Value = 12603
$ws4.Range("F12").Value = 256
$ws4.Range("F13").Value = 7
$ws4.Range("F14").Value = 4936
$ws4.Range("F15").Value = 5621
$ws4.Range("F16").Value = 171
$ws4.Range("F17").Value = 80
$ws4.Range("F20").Value = 974
$ws4.Range("F24").Value = 179
